# Apply the cryptos-list refresh described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") holds plain text in the source data (it may contain
# thousands-separator dots like "67.277.60", or small decimals like
# "0.0000252" that Excel would otherwise render in scientific notation).
# Force the whole column to Text format before writing so COM keeps every
# value as a literal string instead of silently coercing it to a Double,
# then drop the format back to Normal so no stray styling is left behind.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "67.277.60"
$ws.Range("E2").Value = "  +2.21%  "

$ws.Range("D3").Value = "3.108.51"
$ws.Range("E3").Value = "  +5.14%  "

$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").Value = "584.01"
$ws.Range("E5").Value = "  +2.83%  "

$ws.Range("D6").Value = "169.88"
$ws.Range("E6").Value = "  +5.78%  "

$ws.Range("E7").Value = "  -0.06%  "

$ws.Range("D8").Value = "3.107.95"
$ws.Range("E8").Value = "  +5.30%  "

$ws.Range("E9").Value = "  +1.30%  "

$ws.Range("D10").Value = "6.71"
$ws.Range("E10").Value = "  -0.14%  "

$ws.Range("E11").Value = "  +2.28%  "

$ws.Range("E12").Value = "  +6.32%  "

$ws.Range("D13").Value = "0.0000252"
$ws.Range("E13").Value = "  +2.54%  "

$ws.Range("D14").Value = "37.13"
$ws.Range("E14").Value = "  +8.27%  "

$ws.Range("E15").Value = "  -0.17%  "

$ws.Range("D16").Value = "3.617.23"
$ws.Range("E16").Value = "  +4.96%  "

$ws.Range("D17").Value = "67.173.70"
$ws.Range("E17").Value = "  +2.00%  "

$ws.Range("D18").Value = "7.22"
$ws.Range("E18").Value = "  +4.31%  "

$ws.Range("D19").Value = "3.100.60"
$ws.Range("E19").Value = "  +4.98%  "

$ws.Range("D20").Value = "16.15"
$ws.Range("E20").Value = "  +17.54%  "

$ws.Range("D21").Value = "471.58"
$ws.Range("E21").Value = "  +5.82%  "

$ws.Range("D22").Value = "0.716"
$ws.Range("E22").Value = "  +6.16%  "

$ws.Range("D23").Value = "7.54"
$ws.Range("E23").Value = "  +5.27%  "

$ws.Range("D24").Value = "83.62"
$ws.Range("E24").Value = "  +1.73%  "

$ws.Range("B25").Value = "Fetch.AI"
$ws.Range("C25").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D25").Value = "2.32"
$ws.Range("E25").Value = "  +5.37%  "

$ws.Range("B26").Value = "InternetComputer(DFINITY)"
$ws.Range("C26").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D26").Value = "12.88"
$ws.Range("E26").Value = "  +5.73%  "

$ws.Range("D27").Value = "10.22"
$ws.Range("E27").Value = "  +2.00%  "

$ws.Range("E28").Value = "  +0.04%  "

$ws.Range("D29").Value = "8.19"
$ws.Range("E29").Value = "  +1.37%  "

$ws.Range("D30").Value = "2.43"
$ws.Range("E30").Value = "  +4.29%  "

$ws.Range("E31").Value = "  +4.81%  "

$ws.Range("E32").Value = "  +4.04%  "

$ws.Range("D33").Value = "28.66"
$ws.Range("E33").Value = "  +5.74%  "

$ws.Range("E34").Value = "  +5.56%  "

$ws.Range("D35").Value = "0.997"
$ws.Range("E35").Value = "  -0.08%  "

$ws.Range("E36").Value = "  +3.72%  "

$ws.Range("D37").Value = "5.94"
$ws.Range("E37").Value = "  +4.36%  "

$ws.Range("D38").Value = "47.94"
$ws.Range("E38").Value = "  +11.64%  "

$ws.Range("D39").Value = "2.10"
$ws.Range("E39").Value = "  +6.56%  "

$ws.Range("D40").Value = "50.51"
$ws.Range("E40").Value = "  +2.73%  "

$ws.Range("D41").Value = "0.316"
$ws.Range("E41").Value = "  +4.90%  "

$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").Value = "0.123"
$ws.Range("E42").Value = "  +3.70%  "

$ws.Range("B43").Value = "dogwifhat"
$ws.Range("C43").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D43").Value = "2.91"
$ws.Range("E43").Value = "  +3.24%  "

$ws.Range("D44").Value = "8.73"
$ws.Range("E44").Value = "  +4.17%  "

$ws.Range("D45").Value = "393.17"
$ws.Range("E45").Value = "  +2.34%  "

$ws.Range("E46").Value = "  +3.44%  "

$ws.Range("D47").Value = "2.776.51"
$ws.Range("E47").Value = "  +2.36%  "

$ws.Range("D48").Value = "135.26"
$ws.Range("E48").Value = "  +3.73%  "

$ws.Range("E49").Value = "  +0.05%  "

$ws.Range("D50").Value = "24.89"
$ws.Range("E50").Value = "  +7.53%  "

$ws.Range("E51").Value = "  +5.29%  "

# Restore the default style on the price column (keeps the values as text
# without leaving the temporary Text number-format applied).
$priceRange.Style = "Normal"

